# Insert a new data row before row 412 (shifts existing rows 412:495 down to 413:496)
# and populate it with a new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("412:412").Insert()

$newRow = 412

$ws.Cells.Item($newRow, 1).Value = 4
$ws.Cells.Item($newRow, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value = 45244
$ws.Cells.Item($newRow, 5).Value = 10
$ws.Cells.Item($newRow, 6).Value = 100112032
$ws.Cells.Item($newRow, 7).Value = "Zapallo italiano"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 250
$ws.Cells.Item($newRow, 11).Value = 24000
$ws.Cells.Item($newRow, 12).Value = 24000
$ws.Cells.Item($newRow, 13).Value = 24000
$ws.Cells.Item($newRow, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item($newRow, 15).Value = "Región de O'Higgins"
$ws.Cells.Item($newRow, 16).Value = 480
$ws.Cells.Item($newRow, 17).Value = 50
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
